$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 304
$ws.Range("F7").Value = 5422
$ws.Range("F9").Value = 7308
$ws.Range("F11").Value = 50
$ws.Range("F12").Value = 3734
$ws.Range("F13").Value = 60
$ws.Range("F14").Value = 15
$ws.Range("F16").Value = 188
$ws.Range("F17").Value = 135
$ws.Range("F19").Value = 51
$ws.Range("F20").Value = 92
$ws.Range("F22").Value = 3819
$ws.Range("F23").Value = 121
$ws.Range("F24").Value = 5094
$ws.Range("F25").Value = 432
$ws.Range("F26").Value = 2038
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 320
$ws.Range("F29").Value = 7529
$ws.Range("F32").Value = 2125
$ws.Range("F34").Value = 142
$ws.Range("F35").Value = 1152
$ws.Range("F38").Value = 248
$ws.Range("F39").Value = 237
$ws.Range("F41").Value = 1169
$ws.Range("F42").Value = 1168
$ws.Range("F43").Value = 21
$ws.Range("F44").Value = 164
$ws.Range("F45").Value = 1286
$ws.Range("F46").Value = 1972
$ws.Range("F47").Value = 105
$ws.Range("F48").Value = 192
$ws.Range("F49").Value = 1201

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 143
$ws.Range("F9").Value = 933
$ws.Range("F11").Value = 113

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 533
$ws.Range("F3").Value = 711
$ws.Range("F4").Value = 60

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 533
$ws.Range("F6").Value = 711
$ws.Range("F7").Value = 60
$ws.Range("F8").Value = 304
$ws.Range("F9").Value = 5422
$ws.Range("F10").Value = 3734
$ws.Range("F11").Value = 60
$ws.Range("F12").Value = 15
$ws.Range("F14").Value = 188
$ws.Range("F15").Value = 135
$ws.Range("F16").Value = 51
$ws.Range("F17").Value = 92
$ws.Range("F19").Value = 143
$ws.Range("F21").Value = 3819
$ws.Range("F23").Value = 121
$ws.Range("F24").Value = 5094
$ws.Range("F25").Value = 432
$ws.Range("F26").Value = 2038
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 320
$ws.Range("F29").Value = 7530
$ws.Range("F32").Value = 2125
$ws.Range("F34").Value = 142
$ws.Range("F35").Value = 1152
$ws.Range("F37").Value = 248
$ws.Range("F38").Value = 237
$ws.Range("F39").Value = 1169
$ws.Range("F40").Value = 1168
$ws.Range("F41").Value = 21
$ws.Range("F42").Value = 164
$ws.Range("F44").Value = 1286
$ws.Range("F46").Value = 1973
$ws.Range("F47").Value = 105
$ws.Range("F49").Value = 192
